$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet: "ActivityConf" -> "Activity"
$ws.Name = "Activity"

# --- New columns P:T added for the "Task" struct, which spans multiple columns ---

# Row 1: English field names (header)
$ws.Range("P1").Value = "TaskType"
$ws.Range("Q1").Value = "TaskParam1"
$ws.Range("R1").Value = "TaskParam2"
$ws.Range("S1").Value = "TaskParam3"
$ws.Range("T1").Value = "TaskTarget"

# Row 2: field type descriptors
$ws.Range("P2").Value = "{Task}int32"
$ws.Range("Q2").Value = "int32"
$ws.Range("R2").Value = "int32"
$ws.Range("S2").Value = "int32"
$ws.Range("T2").Value = "int32"

# Row 3: Chinese field names
$ws.Range("P3").Value = "任务类型"
$ws.Range("Q3").Value = "参数1"
$ws.Range("R3").Value = "参数2"
$ws.Range("S3").Value = "参数3"
$ws.Range("T3").Value = "目标"

# Row 4: data
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1
$ws.Range("T4").Value = 1

# Row 5: data
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 2
$ws.Range("T5").Value = 1

# Row 6: data
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 3
$ws.Range("T6").Value = 1

# Row 7: no task data for this row

# Row 8: data
$ws.Range("P8").Value = 2
$ws.Range("Q8").Value = 1

# Row 9: data
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 2

# Row 10: data
$ws.Range("P10").Value = 2
$ws.Range("Q10").Value = 3

# Move the active selection, matching the edited file's cursor position
$ws.Range("U9").Select()
